$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header values P1 and Q1 (style matches existing header row, style index 1)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# Update data rows 2-25: columns I, K, M, O change values; add new P, Q columns
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I
    $ws.Cells.Item($r, 11).Value = 1   # K
    $ws.Cells.Item($r, 13).Value = 2   # M
    $ws.Cells.Item($r, 15).Value = 1   # O
    $ws.Cells.Item($r, 16).Value = 2   # P
    $ws.Cells.Item($r, 17).Value = 2   # Q
}
